$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Timer Value gets decremented."
$ws.Range("B9").Value = "Timer Value is now 3."
$ws.Range("A10").Value = "runCountdownTimer() checks if timer value is between 1 - 5."
$ws.Range("B10").Value = "True."
$ws.Range("A11").Value = "runCountdownTimer() checks if audio is on and transition is false."
$ws.Range("B11").Value = "False (because transition is true)."
$ws.Range("A12").Value = "runCountdownTimer() interval is done.  Repeat."
$ws.Range("B12").Value = "timerUI() displays 3 seconds."
$ws.Range("A13").Value = "Timer Value gets decremented."
$ws.Range("B13").Value = "Timer Value is now 2."
$ws.Range("A14").Value = "runCountdownTimer() checks if timer value is between 1 - 5."
$ws.Range("B14").Value = "True."
$ws.Range("A15").Value = "runCountdownTimer() checks if audio is on and transition is false."
$ws.Range("B15").Value = "False (because transition is true)."
$ws.Range("A16").Value = "runCountdownTimer() interval is done.  Repeat."
$ws.Range("B16").Value = "timerUI() displays 2 seconds."
$ws.Range("A17").Value = "Timer Value gets decremented."
$ws.Range("B17").Value = "Timer Value is now 1."
$ws.Range("A18").Value = "runCountdownTimer() checks if timer value is between 1 - 5."
$ws.Range("B18").Value = "True."
$ws.Range("A19").Value = "runCountdownTimer() checks if audio is on and transition is false."
$ws.Range("B19").Value = "False (because transition is true)."
$ws.Range("A20").Value = "runCountdownTimer() interval is done.  Repeat."
$ws.Range("B20").Value = "timerUI() displays 1 seconds."
$ws.Range("A21").Value = "Timer Value gets decremented."
$ws.Range("B21").Value = "Timer Value is now 0."
$ws.Range("A22").Value = "runCountdownTimer() checks if timer value is between 1 - 5."
$ws.Range("B22").Value = "False (because timer value is 0)."
$ws.Range("A23").Value = "runCountdownTimer() checks if audio is on."
$ws.Range("B23").Value = "If true, plays audio.  If false, does not."
$ws.Range("A24").Value = "runCountdownTimer() checks if transition is true."
$ws.Range("B24").Value = "True."
$ws.Range("A25").Value = "runCountdownTimer() sets transition to false."
$ws.Range("B25").Value = "transition is now false."
$ws.Range("A26").Value = "runCountdownTimer() interval is done.  Repeat."

# The three trailing rows in the source workbook carry an explicit
# (visually no-op) "no fill" cell style on column B - reproduce that so the
# new cellXfs entry / applyFill flag shows up the same way it did upstream.
# (White keeps the cells visually blank/unfilled against the sheet's default
# white background, matching the original "no fill" look.)
$ws.Range("B23:B25").Interior.Color = 16777215

$null = $ws.Range("A26").Select()
